$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite1")

# ---------------------------------------------------------------------------
# 1. Insert a new column K ("Parameters to write") before the existing
#    "Actual Result" / "Pass/Fail" columns, shifting them right.
# ---------------------------------------------------------------------------
$ws.Range("K1").EntireColumn.Insert()
$ws.Range("K1").Value = "Parameters to write"
$ws.Range("K2").Clear()

# ---------------------------------------------------------------------------
# 2. Row 2 (TC01) content updates.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "P_Valid  Login"
$ws.Range("H2").Value = "<?xml version=""1.0"" encoding=""UTF-8""?>
<SOAP-ENV:Envelope xmlns:SOAP-ENV=""http://schemas.xmlsoap.org/soap/envelope/"" xmlns:ns1=""https://api.shmart.in"">
    <SOAP-ENV:Body>
        <ns1:LoginResponse>
            <return>
                <SessionID>B6B1AA1A99</SessionID>
                <ResponseCode>0</ResponseCode>
                <ResponseMessage>Successful</ResponseMessage>
            </return>
        </ns1:LoginResponse>
    </SOAP-ENV:Body>
</SOAP-ENV:Envelope>"
$ws.Range("J2").Value = "SessionID"
$ws.Rows("2").RowHeight = 168.75

# ---------------------------------------------------------------------------
# 3. Insert new row 3 (TC02 - N_invalid_Login) below row 2.
# ---------------------------------------------------------------------------
$ws.Rows("3").Insert()

$ws.Range("A5").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("G3").WrapText = $true
$ws.Range("H3").WrapText = $true
$ws.Range("I3").Clear()
$ws.Range("K3").Clear()

$ws.Range("A3").Value = "TC02"
$ws.Range("B3").Value = "Login"
$ws.Range("C3").Value = "N_invalid_Login"
$ws.Range("D3").Value = "N_invalid_Login"
$ws.Range("E3").Value = "SOAP"
$ws.Range("F3").Value = "/services/partners"
$ws.Range("G3").Value = "<SOAP-ENV:Envelope xmlns:SOAP-ENV=""http://schemas.xmlsoap.org/soap/envelope/""
xmlns:sas=""http://api.shmart.in/"">
<SOAP-ENV:Body>
        <sas:Login>
                <Username>#random_string</Username>
                <Password>pratik</Password>
        </sas:Login>
</SOAP-ENV:Body>
</SOAP-ENV:Envelope>
"
$ws.Range("H3").Value = "<?xml version=""1.0"" encoding=""UTF-8""?>
<SOAP-ENV:Envelope xmlns:SOAP-ENV=""http://schemas.xmlsoap.org/soap/envelope/"" xmlns:ns1=""https://api.shmart.in"">
    <SOAP-ENV:Body>
        <ns1:LoginResponse>
            <return>
                <ResponseMessage>Invalid Login</ResponseMessage>
                <ResponseCode>100</ResponseCode>
            </return>
        </ns1:LoginResponse>
    </SOAP-ENV:Body>
</SOAP-ENV:Envelope>"
$ws.Range("J3").Value = "null"
$ws.Rows("3").RowHeight = 146.25

# ---------------------------------------------------------------------------
# 4. Insert new row 4 (TC03 - N_invalid_Session) below row 3.
# ---------------------------------------------------------------------------
$ws.Rows("4").Insert()

$ws.Range("A5").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("G4").WrapText = $true
$ws.Range("H4").WrapText = $true
$ws.Range("I4").Clear()
$ws.Range("K4").Clear()

$ws.Range("A4").Value = "TC03"
$ws.Range("B4").Value = "Login"
$ws.Range("C4").Value = "N_invalid_Session"
$ws.Range("D4").Value = "N_invalid_Session"
$ws.Range("E4").Value = "SOAP"
$ws.Range("F4").Value = "/services/partners"
$ws.Range("G4").Value = "<SOAP-ENV:Envelope xmlns:SOAP-ENV=""http://schemas.xmlsoap.org/soap/envelope/""
xmlns:sas=""http://api.shmart.in/"">
<SOAP-ENV:Body>
        <sas:GenerateOTPRequest>
                <SessionID>DC421FA156</SessionID>
                <ProductCode>27</ProductCode>
                <Mobile>9833868977</Mobile>
                <RequestType>R</RequestType>
         </sas:GenerateOTPRequest>
</SOAP-ENV:Body>
</SOAP-ENV:Envelope>
"
$ws.Range("H4").Value = "<?xml version=""1.0"" encoding=""UTF-8""?>
<SOAP-ENV:Envelope xmlns:SOAP-ENV=""http://schemas.xmlsoap.org/soap/envelope/"" xmlns:ns1=""https://api.shmart.in"">
    <SOAP-ENV:Body>
        <ns1:GenerateOTPRequestResponse>
            <return>
                <ResponseMessage>Invalid Login</ResponseMessage>
                <ResponseCode>100</ResponseCode>
            </return>
        </ns1:GenerateOTPRequestResponse>
    </SOAP-ENV:Body>
</SOAP-ENV:Envelope>"
$ws.Range("J4").Value = "null"
$ws.Rows("4").RowHeight = 157.5

# ---------------------------------------------------------------------------
# 5. Column width + view adjustments.
# ---------------------------------------------------------------------------
$ws.Range("H1").ColumnWidth = 43.6

$excel.ActiveWindow.ScrollColumn = 3
[void]$ws.Range("H3").Select()

Write-Output "done"
